$wb = $excel.ActiveWorkbook

# --- "BoM" sheet: reorder comma separated Net Name lists in column X ---
$wsBom = $wb.Worksheets.Item("BoM")

$wsBom.Range("X9").Value  = "+5V,GND"
$wsBom.Range("X11").Value = "+5V,GND"
$wsBom.Range("X13").Value = "/RXLED,Net-(D2-A)"
$wsBom.Range("X14").Value = ",+5V,GND"
$wsBom.Range("X15").Value = "+5V,GND,/SCK2,/MOSI2,/MISO2,/RESET2"
$wsBom.Range("X16").Value = "Net-(J3-Pin_1),Net-(J3-Pin_5),Net-(J3-Pin_2),Net-(J3-Pin_3),Net-(J3-Pin_4)"
$wsBom.Range("X17").Value = "Net-(J6-Pin_5),Net-(J6-Pin_2),Net-(J6-Pin_6),Net-(J6-Pin_3),Net-(J6-Pin_4),Net-(J6-Pin_1)"
$wsBom.Range("X21").Value = "Net-(J4-Pin_2),Net-(U1-D+),GND,/DTR,unconnected-(U1-PB0-Pad14),Net-(J6-Pin_3),Net-(J3-Pin_3),Net-(U1-PC0{slash}XTAL2),Net-(J6-Pin_5),Net-(J6-Pin_2),Net-(J6-Pin_6),Net-(U1-UCAP),/SCK2,/MOSI2,Net-(J3-Pin_2),/RESET2,Earth,Net-(J4-Pin_1),/MISO2,Net-(J4-Pin_3),Net-(J3-Pin_4),+5V,Net-(J3-Pin_1),Net-(J3-Pin_5),VBUS,/RXLED,Net-(J4-Pin_4),Net-(J6-Pin_4),/TXLED,Net-(U1-D-),Net-(U1-XTAL1)"

# Row 21 is taller now to fit the (still long) rewrapped net list
$wsBom.Rows.Item(21).RowHeight = 120

# --- "DNF" sheet: reorder comma separated Net Name lists in column X ---
$wsDnf = $wb.Worksheets.Item("DNF")

$wsDnf.Range("X9").Value  = "GND,Net-(U1-XTAL1)"
$wsDnf.Range("X12").Value = "Net-(J4-Pin_2),Net-(J4-Pin_4),Net-(J4-Pin_3),Net-(J4-Pin_1)"
$wsDnf.Range("X13").Value = "Net-(J2-VBUS),Net-(J2-D+),Earth,Net-(J2-D-),Net-(J2-Shield)"
$wsDnf.Range("X14").Value = "Net-(U1-D+),Net-(J2-D+)"
$wsDnf.Range("X15").Value = "Net-(U1-PC0{slash}XTAL2),Net-(U1-XTAL1)"
$wsDnf.Range("X16").Value = "Net-(J2-Shield),Net-(J2-D+)"
$wsDnf.Range("X17").Value = "Net-(U1-PC0{slash}XTAL2),Net-(U1-XTAL1)"
